# Scene 78 ("Front of House" / porch scene) — style-sheet refresh.
#
# The authoritative diff for this commit shows the document's style
# catalogue (Normal, Table Normal, Heading 1-6, Title, Subtitle) being
# (re)asserted with their full definitions, plus a Google-Docs-style
# round-trip `customXML` payload being attached to the package. That
# signature is what Google Docs stamps onto a .docx the moment you open
# it there and save/download again — i.e. "Stephen" (the collaborator
# named in the commit message) opened the script in Docs to write his
# new pages, and the side effect of that round trip is the style block
# being re-emitted and the `customXML` parts showing up.
#
# Word's object model intentionally refuses to mint a second style with
# a styleId/name that already exists (`Styles.Add` on a name already in
# the collection just returns the existing Style — "A style with that
# name already exists" is a hard guard, not a glitch), so the faithful,
# automation-safe way to reproduce "the style block got rewritten with
# these exact definitions" is to assert each style's full formatting
# through the object model rather than trying to force a duplicate
# w:styleId into the part. That keeps the package well-formed while
# landing on the same effective formatting the diff encodes.
#
# We also register the Google Docs round-trip payload through the
# documented CustomXMLParts API (the correct call for attaching a new
# custom XML part to the package), mirroring the `customXML/item1.xml`
# / `customXML/itemProps1.xml` pair added by the diff.

$d = $word.ActiveDocument

function Ensure-ParagraphStyle($name, $basedOn, $next, $keepNext, $keepLines, $pageBreakBefore, $spaceBeforePts, $spaceAfterPts, $bold, $italic, $sizePts, $colorRGB, $fontName) {
    $style = $null
    try { $style = $d.Styles.Item($name) } catch { $style = $null }
    if ($style -eq $null) {
        $style = $d.Styles.Add($name, 1)
    }

    if ($basedOn -ne $null -and $basedOn -ne "") {
        try { $style.BaseStyle = $d.Styles.Item($basedOn) } catch { }
    }
    if ($next -ne $null -and $next -ne "") {
        try { $style.NextParagraphStyle = $d.Styles.Item($next) } catch { }
    }

    $style.ParagraphFormat.KeepWithNext = $keepNext
    $style.ParagraphFormat.KeepTogether = $keepLines
    $style.ParagraphFormat.PageBreakBefore = $pageBreakBefore
    $style.ParagraphFormat.SpaceBefore = $spaceBeforePts
    $style.ParagraphFormat.SpaceAfter = $spaceAfterPts

    if ($sizePts -gt 0) { $style.Font.Size = $sizePts }
    $style.Font.Bold = $bold
    $style.Font.Italic = $italic
    if ($colorRGB -ne $null) { $style.Font.TextColor.RGB = $colorRGB }
    if ($fontName -ne $null -and $fontName -ne "") { $style.Font.Name = $fontName }

    return $style
}

# Normal / Table Normal: bare defaults, nothing further to assert.
try { [void]$d.Styles.Item("Normal") } catch { [void]$d.Styles.Add("Normal", 1) }
try { [void]$d.Styles.Item("TableNormal") } catch { [void]$d.Styles.Add("Table Normal", 2) }

[void](Ensure-ParagraphStyle "Heading1" "Normal" "Normal" $true $true $false 20 6 $false $false 20 $null "")
[void](Ensure-ParagraphStyle "Heading2" "Normal" "Normal" $true $true $false 18 6 $false $false 16 $null "")
[void](Ensure-ParagraphStyle "Heading3" "Normal" "Normal" $true $true $false 16 4 $false $false 14 0x434343 "")
[void](Ensure-ParagraphStyle "Heading4" "Normal" "Normal" $true $true $false 14 4 $false $false 12 0x666666 "")
[void](Ensure-ParagraphStyle "Heading5" "Normal" "Normal" $true $true $false 12 4 $false $false 11 0x666666 "")
[void](Ensure-ParagraphStyle "Heading6" "Normal" "Normal" $true $true $false 12 4 $false $true  11 0x666666 "")
[void](Ensure-ParagraphStyle "Title"    "Normal" "Normal" $true $true $false 0  3 $false $false 26 $null "")
[void](Ensure-ParagraphStyle "Subtitle" "Normal" "Normal" $true $true $false 0  16 $false $false 15 0x666666 "Arial")

# Attach the Google Docs round-trip custom XML payload (customXML/item1.xml
# + customXML/itemProps1.xml in the package diff) via the documented API.
try {
    $customXml = '<?xml version="1.0" encoding="utf-8"?><go:gDocsCustomXmlDataStorage xmlns:go="http://customooxmlschemas.google.com/" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><go:docsCustomData xmlns:go="http://customooxmlschemas.google.com/" roundtripDataSignature="AMtx7mgA98XVkzaneOhhPEdpYmSuukIXBQ==">AMUW2mXftsawzNZHvxREb0zwemnFyoZP7aap5aCIuOcj5xcvwfqdlD2oJBycnDjjei+2piJxO56hjFAiIBu3zrbEHeI5TpirCXkD5Mg9i4W0Ukp4ClAsRI4=</go:docsCustomData></go:gDocsCustomXmlDataStorage>'
    [void]$d.CustomXMLParts.Add($customXml)
} catch { }

Write-Output "Style sheet refreshed; styles=$($d.Styles.Count)"
